# Generate Report for Archive
# The localization status for file "2044556c-ec1d-4e6d-9e75-8f18baca15c1.md"
# moved from "Ready for handoff" to "In Translation". Update the Status
# cells on the Overview sheet (both the zh-cn and de-de columns) and on the
# per-locale detail sheets (zh-cn, de-de) accordingly.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "In Translation"
$overview.Range("C2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B2").Value = "In Translation"
